# Auto-generated edit script
# Applies updated market/profit values to the Tiamat_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3085.8572
$ws.Range("I70").Value = 3050
$ws.Range("J70").Value = 3112.75
$ws.Range("K70").Value = 9150
$ws.Range("L70").Value = 9338.25
$ws.Range("M70").Value = -8880
$ws.Range("N70").Value = -9878.25
$ws.Range("H73").Value = 3085.8572
$ws.Range("I73").Value = 3050
$ws.Range("J73").Value = 3112.75
$ws.Range("K73").Value = 9150
$ws.Range("L73").Value = 9338.25
$ws.Range("M73").Value = -8214
$ws.Range("N73").Value = -11210.25
$ws.Range("H134").Value = 60000
$ws.Range("J134").Value = 60000
$ws.Range("L134").Value = 60000
$ws.Range("N134").Value = -70140
$ws.Range("H137").Value = 3595.0256
$ws.Range("I137").Value = 820.7692
$ws.Range("J137").Value = 9143.538
$ws.Range("K137").Value = 2462.3076
$ws.Range("L137").Value = 27430.614
$ws.Range("M137").Value = 87.69239999999991
$ws.Range("N137").Value = -32530.614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 14934.333
$ws.Range("I3").Value = 2427.5
$ws.Range("J3").Value = 39948
$ws.Range("K3").Value = 2427.5
$ws.Range("L3").Value = 39948
$ws.Range("M3").Value = -2312.5
$ws.Range("N3").Value = -40178
$ws.Range("H61").Value = 2792.625
$ws.Range("I61").Value = 2961
$ws.Range("J61").Value = 1614
$ws.Range("K61").Value = 2961
$ws.Range("L61").Value = 1614
$ws.Range("M61").Value = -2749
$ws.Range("N61").Value = -2038
$ws.Range("H74").Value = 1058.32
$ws.Range("I74").Value = 1096.9678
$ws.Range("J74").Value = 995.2632
$ws.Range("K74").Value = 1096.9678
$ws.Range("L74").Value = 995.2632
$ws.Range("M74").Value = -222.9677999999999
$ws.Range("N74").Value = -2743.2632
$ws.Range("H77").Value = 1058.32
$ws.Range("I77").Value = 1096.9678
$ws.Range("J77").Value = 995.2632
$ws.Range("K77").Value = 5484.839
$ws.Range("L77").Value = 4976.316
$ws.Range("M77").Value = -1116.839
$ws.Range("N77").Value = -13712.316
$ws.Range("H132").Value = 26159.979
$ws.Range("I132").Value = 35886.793
$ws.Range("J132").Value = 8530.125
$ws.Range("K132").Value = 107660.379
$ws.Range("L132").Value = 25590.375
$ws.Range("M132").Value = -105130.379
$ws.Range("N132").Value = -30650.375
$ws.Range("H136").Value = 2792.625
$ws.Range("I136").Value = 2961
$ws.Range("J136").Value = 1614
$ws.Range("K136").Value = 8883
$ws.Range("L136").Value = 4842
$ws.Range("M136").Value = -6333
$ws.Range("N136").Value = -9942

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1766.6666
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253
$ws.Range("H134").Value = 4187.067
$ws.Range("I134").Value = 2101.2
$ws.Range("J134").Value = 8358.799999999999
$ws.Range("K134").Value = 6303.599999999999
$ws.Range("L134").Value = 25076.4
$ws.Range("M134").Value = -3768.599999999999
$ws.Range("N134").Value = -30146.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 27737.143
$ws.Range("I2").Value = 383.33334
$ws.Range("J2").Value = 48252.5
$ws.Range("K2").Value = 383.33334
$ws.Range("L2").Value = 48252.5
$ws.Range("M2").Value = -270.33334
$ws.Range("N2").Value = -48478.5
$ws.Range("H31").Value = 24289.432
$ws.Range("I31").Value = 24809.76
$ws.Range("K31").Value = 24809.76
$ws.Range("M31").Value = -24514.76
$ws.Range("H34").Value = 24289.432
$ws.Range("I34").Value = 24809.76
$ws.Range("K34").Value = 24809.76
$ws.Range("M34").Value = -24607.76
$ws.Range("H58").Value = 2001.1765
$ws.Range("I58").Value = 1079
$ws.Range("J58").Value = 6304.6665
$ws.Range("K58").Value = 1079
$ws.Range("L58").Value = 6304.6665
$ws.Range("M58").Value = -876
$ws.Range("N58").Value = -6710.6665
$ws.Range("H62").Value = 125003250
$ws.Range("J62").Value = 4500
$ws.Range("L62").Value = 4500
$ws.Range("N62").Value = -5748
$ws.Range("H65").Value = 125003250
$ws.Range("J65").Value = 4500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -28740
$ws.Range("H132").Value = 2054.3235
$ws.Range("I132").Value = 1305.08
$ws.Range("J132").Value = 4135.5557
$ws.Range("K132").Value = 3915.24
$ws.Range("L132").Value = 12406.6671
$ws.Range("M132").Value = -1385.24
$ws.Range("N132").Value = -17466.6671
$ws.Range("H134").Value = 25001964
$ws.Range("I134").Value = 1699.1333
$ws.Range("J134").Value = 100002760
$ws.Range("K134").Value = 5097.3999
$ws.Range("L134").Value = 300008280
$ws.Range("M134").Value = -2562.3999
$ws.Range("N134").Value = -300013350
$ws.Range("H136").Value = 2001.1765
$ws.Range("I136").Value = 1079
$ws.Range("J136").Value = 6304.6665
$ws.Range("K136").Value = 3237
$ws.Range("L136").Value = 18913.9995
$ws.Range("M136").Value = -687
$ws.Range("N136").Value = -24013.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 45455996
$ws.Range("I44").Value = 331.66666
$ws.Range("J44").Value = 62501868
$ws.Range("K44").Value = 994.9999799999999
$ws.Range("L44").Value = 187505604
$ws.Range("M44").Value = -596.9999799999999
$ws.Range("N44").Value = -187506400
$ws.Range("H64").Value = 4334162
$ws.Range("J64").Value = 5056105.5
$ws.Range("L64").Value = 15168316.5
$ws.Range("N64").Value = -15168856.5
$ws.Range("H67").Value = 4334162
$ws.Range("J67").Value = 5056105.5
$ws.Range("L67").Value = 15168316.5
$ws.Range("N67").Value = -15170188.5
$ws.Range("H107").Value = 255.15384
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 255.15384
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 765.4615200000001
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4605.46152

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 53586.75
$ws.Range("I132").Value = 2925
$ws.Range("J132").Value = 87361.25
$ws.Range("K132").Value = 8775
$ws.Range("L132").Value = 262083.75
$ws.Range("M132").Value = -6245
$ws.Range("N132").Value = -267143.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 27090.635
$ws.Range("I132").Value = 44184.332
$ws.Range("J132").Value = 2958.353
$ws.Range("K132").Value = 132552.996
$ws.Range("L132").Value = 8875.059000000001
$ws.Range("M132").Value = -130022.996
$ws.Range("N132").Value = -13935.059
$ws.Range("H136").Value = 1783.6666
$ws.Range("I136").Value = 963.61536
$ws.Range("J136").Value = 3116.25
$ws.Range("K136").Value = 2890.84608
$ws.Range("L136").Value = 9348.75
$ws.Range("M136").Value = -340.8460800000003
$ws.Range("N136").Value = -14448.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1618.675
$ws.Range("I132").Value = 1470.9565
$ws.Range("J132").Value = 1818.5294
$ws.Range("K132").Value = 4412.8695
$ws.Range("L132").Value = 5455.5882
$ws.Range("M132").Value = -1882.8695
$ws.Range("N132").Value = -10515.5882
$ws.Range("H136").Value = 3132926
$ws.Range("I136").Value = 3402524.8
$ws.Range("J136").Value = 2000611.8
$ws.Range("K136").Value = 10207574.4
$ws.Range("L136").Value = 6001835.4
$ws.Range("M136").Value = -10205024.4
$ws.Range("N136").Value = -6006935.4
